$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BAC-18")

# Update the content: "Contento de trabajadores" -> "Rotación de trabajadores"
$ws.Range("B8").Value = "Rotación de trabajadores"

# Update the active selection to match the saved view (B9)
$ws.Range("B9").Select()
